# Scheduled market-data refresh: update cached price/profit figures on the
# Tiamat_Profits sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 2312.7144
$ws.Range("I29").Value = 2243.8
$ws.Range("J29").Value = 2485
$ws.Range("K29").Value = 6731.400000000001
$ws.Range("L29").Value = 7455
$ws.Range("M29").Value = -6450.400000000001
$ws.Range("N29").Value = -8017
# Row 105
$ws.Range("H105").Value = 28000
$ws.Range("J105").Value = 28000
$ws.Range("L105").Value = 28000
$ws.Range("N105").Value = -34988
# Row 116
$ws.Range("H116").Value = 7448.048
$ws.Range("I116").Value = 16429
$ws.Range("K116").Value = 16429
$ws.Range("M116").Value = -12987
# Row 129
$ws.Range("H129").Value = 729.55554
$ws.Range("J129").Value = 1025.7142
$ws.Range("L129").Value = 3077.1426
$ws.Range("N129").Value = -13077.1426
# Row 138
$ws.Range("H138").Value = 8791715
$ws.Range("I138").Value = 2407.818
$ws.Range("J138").Value = 10893506
$ws.Range("K138").Value = 7223.454000000001
$ws.Range("L138").Value = 32680518
$ws.Range("M138").Value = -2083.454000000001
$ws.Range("N138").Value = -32690798

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 883.13043
$ws.Range("I61").Value = 829.1429000000001
$ws.Range("K61").Value = 829.1429000000001
$ws.Range("M61").Value = -617.1429000000001
# Row 97
$ws.Range("H97").Value = 3111
$ws.Range("I97").Value = 2166.5
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 2166.5
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -1670.5
$ws.Range("N97").Value = -5992
# Row 136
$ws.Range("H136").Value = 883.13043
$ws.Range("I136").Value = 829.1429000000001
$ws.Range("K136").Value = 2487.4287
$ws.Range("M136").Value = 62.57129999999961

$ws = $wb.Worksheets.Item("CRP")
# Row 74
$ws.Range("H74").Value = 11571.385
$ws.Range("J74").Value = 11571.385
$ws.Range("L74").Value = 11571.385
$ws.Range("N74").Value = -13319.385
# Row 77
$ws.Range("H77").Value = 11571.385
$ws.Range("J77").Value = 11571.385
$ws.Range("L77").Value = 34714.155
$ws.Range("N77").Value = -43450.155
# Row 92
$ws.Range("H92").Value = 24666.666
$ws.Range("J92").Value = 24666.666
$ws.Range("L92").Value = 24666.666
$ws.Range("N92").Value = -29658.666
# Row 105
$ws.Range("H105").Value = 744.7143
$ws.Range("I105").Value = 764.6667
$ws.Range("J105").Value = 625
$ws.Range("K105").Value = 764.6667
$ws.Range("L105").Value = 625
$ws.Range("M105").Value = 982.3333
$ws.Range("N105").Value = -4119
# Row 134
$ws.Range("H134").Value = 854.9322
$ws.Range("I134").Value = 847.05554
$ws.Range("J134").Value = 940
$ws.Range("K134").Value = 2541.16662
$ws.Range("L134").Value = 2820
$ws.Range("M134").Value = -6.166619999999966
$ws.Range("N134").Value = -7890

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 8545.691999999999
$ws.Range("I5").Value = 899.1429000000001
$ws.Range("J5").Value = 17466.666
$ws.Range("K5").Value = 2697.4287
$ws.Range("L5").Value = 52399.99800000001
$ws.Range("M5").Value = -2585.4287
$ws.Range("N5").Value = -52623.99800000001
# Row 33
$ws.Range("H33").Value = 2164.7646
$ws.Range("I33").Value = 696.55554
$ws.Range("J33").Value = 3816.5
$ws.Range("K33").Value = 4179.33324
$ws.Range("L33").Value = 22899
$ws.Range("M33").Value = -3896.33324
$ws.Range("N33").Value = -23465
# Row 116
$ws.Range("H116").Value = 2492.8572
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
# Row 131
$ws.Range("H131").Value = 39063332
$ws.Range("I131").Value = 552.8570999999999
$ws.Range("J131").Value = 50000908
$ws.Range("K131").Value = 1658.5713
$ws.Range("L131").Value = 150002724
$ws.Range("M131").Value = 3381.4287
$ws.Range("N131").Value = -150012804
# Row 135
$ws.Range("H135").Value = 8545.691999999999
$ws.Range("I135").Value = 899.1429000000001
$ws.Range("J135").Value = 17466.666
$ws.Range("K135").Value = 8092.2861
$ws.Range("L135").Value = 157199.994
$ws.Range("M135").Value = -5557.2861
$ws.Range("N135").Value = -162269.994

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 11275
$ws.Range("I43").Value = 1217.4
$ws.Range("J43").Value = 19656.334
$ws.Range("K43").Value = 1217.4
$ws.Range("L43").Value = 19656.334
$ws.Range("M43").Value = -1066.4
$ws.Range("N43").Value = -19958.334
# Row 46
$ws.Range("H46").Value = 11250
$ws.Range("J46").Value = 11250
$ws.Range("L46").Value = 11250
$ws.Range("N46").Value = -11562
# Row 57
$ws.Range("H57").Value = 9333.333000000001
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 80
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2002
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10008
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 575.9167
$ws.Range("I46").Value = 450.16666
$ws.Range("J46").Value = 701.6667
$ws.Range("K46").Value = 450.16666
$ws.Range("L46").Value = 701.6667
$ws.Range("M46").Value = -262.16666
$ws.Range("N46").Value = -1077.6667
# Row 68
$ws.Range("H68").Value = 2282.353
$ws.Range("I68").Value = 1743.75
$ws.Range("J68").Value = 2761.111
$ws.Range("K68").Value = 1743.75
$ws.Range("L68").Value = 2761.111
$ws.Range("M68").Value = -994.75
$ws.Range("N68").Value = -4259.111
# Row 71
$ws.Range("H71").Value = 2282.353
$ws.Range("I71").Value = 1743.75
$ws.Range("J71").Value = 2761.111
$ws.Range("K71").Value = 8718.75
$ws.Range("L71").Value = 13805.555
$ws.Range("M71").Value = -4974.75
$ws.Range("N71").Value = -21293.555
# Row 82
$ws.Range("H82").Value = 1827.238
$ws.Range("I82").Value = 1314
$ws.Range("K82").Value = 1314
$ws.Range("M82").Value = -953
# Row 85
$ws.Range("H85").Value = 1827.238
$ws.Range("I85").Value = 1314
$ws.Range("K85").Value = 1314
$ws.Range("M85").Value = -66
# Row 136
$ws.Range("H136").Value = 669097.5600000001
$ws.Range("I136").Value = 1251317.9
$ws.Range("J136").Value = 3702.8572
$ws.Range("K136").Value = 3753953.7
$ws.Range("L136").Value = 11108.5716
$ws.Range("M136").Value = -3751403.7
$ws.Range("N136").Value = -16208.5716

$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 43002.2
$ws.Range("I8").Value = 32501.5
$ws.Range("J8").Value = 50002.668
$ws.Range("K8").Value = 32501.5
$ws.Range("L8").Value = 50002.668
$ws.Range("M8").Value = -32361.5
$ws.Range("N8").Value = -50282.668
# Row 10
$ws.Range("H10").Value = 5000
$ws.Range("J10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("N10").Value = -5338
